# UTS - Perbaikan error/bug
# Keep only one data row (the "staff" user), renumber its level_id to 3,
# update its username/nama, and remove the now-obsolete rows 3 and 4.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 in place.
# Note: set C2 before B2 so new shared strings are appended in the same
# order as the target workbook (bella-Pontianak, then staff-84).
$ws.Range("A2").Value = 3
$ws.Range("C2").Value = "bella-Pontianak"
$ws.Range("B2").Value = "staff-84"

# Remove the old rows 3 and 4 entirely (data now fits in a single row).
$ws.Rows("3:4").Delete()

# Match the saved selection/active cell state.
$ws.Range("D7").Select() | Out-Null
